$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append (row 4)
$ws.Range("A4").Value = "0W00ChFL"
# B4 holds a date-like literal string ("08/11/2024"); force text formatting
# before assignment so Excel doesn't auto-convert it to a date serial, then
# restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "08/11/2024"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "10:30"
$ws.Range("D4").Value = "UKRAINE - PREMIER LEAGUE"
$ws.Range("E4").Value = "Ch. Odesa"
$ws.Range("F4").Value = "Vorskla Poltava"
$ws.Range("G4").Value = 3.55
$ws.Range("H4").Value = 3.05
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 4.2
$ws.Range("K4").Value = 1.91
$ws.Range("L4").Value = 2.77
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 7.4
$ws.Range("O4").Value = 1.45
$ws.Range("P4").Value = 2.4
$ws.Range("Q4").Value = 2.27
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 1.98
$ws.Range("V4").Value = 1.65
$ws.Range("W4").Value = 8
$ws.Range("X4").Value = 17
$ws.Range("Y4").Value = 13
$ws.Range("Z4").Value = 50
$ws.Range("AA4").Value = 40
$ws.Range("AB4").Value = 55
$ws.Range("AC4").Value = 6.9
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 18
$ws.Range("AF4").Value = 110
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 5.9
$ws.Range("AI4").Value = 9
$ws.Range("AJ4").Value = 9.25
$ws.Range("AK4").Value = 19.5
$ws.Range("AL4").Value = 20
$ws.Range("AM4").Value = 40
$ws.Range("AN4").Value = 5.2
$ws.Range("AO4").Value = 21
$ws.Range("AP4").Value = 32
$ws.Range("AQ4").Value = 120
$ws.Range("AR4").Value = 200
$ws.Range("AS4").Value = 500
$ws.Range("AT4").Value = 2.2
$ws.Range("AU4").Value = 7.9
$ws.Range("AV4").Value = 90
$ws.Range("AW4").Value = 3.75
$ws.Range("AX4").Value = 11.25
$ws.Range("AY4").Value = 24
$ws.Range("AZ4").Value = 50
$ws.Range("BA4").Value = 100
$ws.Range("BB4").Value = 400
$ws.Range("BC4").Value = 81
$ws.Range("BD4").Value = 81

# Update existing cell I3 value from 1.33 to 1.3
$ws.Range("I3").Value = 1.3
